# Operation Chaining 6ns vivado report aggiornati
# Update the DSP power report table (rows 2-6) with refreshed Vivado numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - top-level utilization total changed
$ws.Range("A2").Value = [double]"2.449857711326331E-4"

# Row 3 - now reports buff1_reg (was buff2_reg__0); clock renamed to the BUFG net
$ws.Range("A3").Value = [double]"1.1746238305931911E-4"
$ws.Range("B3").Value = "buff1_reg (DSP48E1)"
$ws.Range("D3").Value = "ap_clk_IBUF_BUFG"
$ws.Range("H3").Value = [double]"5.282081127166748"

# Row 4 - stays buff2_reg__0; clock renamed to the BUFG net
$ws.Range("A4").Value = [double]"1.0035662853624672E-4"
$ws.Range("B4").Value = "buff2_reg__0 (DSP48E1)"
$ws.Range("D4").Value = "ap_clk_IBUF_BUFG"
$ws.Range("H4").Value = [double]"4.465609073638916"

# Row 5 - now reports buff2_reg__0 (was buff1_reg); clock renamed to the BUFG net
$ws.Range("A5").Value = [double]"1.8320855815545656E-5"
$ws.Range("B5").Value = "buff2_reg__0 (DSP48E1)"
$ws.Range("D5").Value = "ap_clk_IBUF_BUFG"
$ws.Range("H5").Value = [double]"0.857142984867096"

# Row 6 - stays buff1_reg; clock renamed to the BUFG net
$ws.Range("A6").Value = [double]"8.84591463545803E-6"
$ws.Range("B6").Value = "buff1_reg (DSP48E1)"
$ws.Range("D6").Value = "ap_clk_IBUF_BUFG"
$ws.Range("H6").Value = [double]"0.42838001251220703"
